# Update "想去人数" (F column) figures across the four worksheets to match
# the freshly generated data snapshot (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

# 展览 (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$updates1 = @{
    2  = 792
    3  = 2861
    4  = 1347
    6  = 595
    11 = 11879
    12 = 6738
    15 = 429
    19 = 937
    23 = 3669
    25 = 992
    26 = 500
    27 = 184
    30 = 234
    32 = 36
    34 = 5053
    35 = 49
    36 = 1259
    37 = 246
    38 = 591
    39 = 218
    40 = 551
    41 = 62
}
foreach ($row in $updates1.Keys) {
    $ws1.Range("F$row").Value = $updates1[$row]
}

# 演出 (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$updates2 = @{
    12 = 3696
    15 = 13
    21 = 5
}
foreach ($row in $updates2.Keys) {
    $ws2.Range("F$row").Value = $updates2[$row]
}

# 本地生活 (Local life)
$ws3 = $wb.Worksheets.Item("本地生活")
$updates3 = @{
    2 = 9100
    4 = 1855
}
foreach ($row in $updates3.Keys) {
    $ws3.Range("F$row").Value = $updates3[$row]
}

# 全部类型 (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$updates4 = @{
    2  = 9100
    4  = 1855
    5  = 792
    6  = 2861
    10 = 1347
    11 = 595
    17 = 11879
    18 = 6738
    19 = 3696
    22 = 429
    26 = 937
    29 = 3669
    31 = 992
    32 = 184
    35 = 234
    39 = 1259
    40 = 246
    41 = 218
    42 = 551
    43 = 5
    46 = 62
}
foreach ($row in $updates4.Keys) {
    $ws4.Range("F$row").Value = $updates4[$row]
}
